$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 415, shifting the
# remaining data (old rows 415:532) down to 418:535.
$ws.Rows("415:417").Insert()

# New row 415: Extra quality record dated 2022-12-05 (serial 44900)
$ws.Range("A415").Value2 = 3
$ws.Range("B415").Value2 = "Femacal de La Calera"
$ws.Range("C415").Value2 = "Coquimbo"
$ws.Range("D415").Value2 = 44900
$ws.Range("E415").Value2 = 5
$ws.Range("F415").Value2 = 100112028
$ws.Range("G415").Value2 = "Sandia"
$ws.Range("H415").Value2 = "Sin especificar"
$ws.Range("I415").Value2 = "Extra"
$ws.Range("J415").Value2 = 550
$ws.Range("K415").Value2 = 4000
$ws.Range("L415").Value2 = 4000
$ws.Range("M415").Value2 = 4000
$ws.Range("N415").Value2 = "`$/unidad"
$ws.Range("O415").Value2 = "Paine"
$ws.Range("P415").Value2 = 4000
$ws.Range("Q415").Value2 = 1
$ws.Range("R415").Value2 = "Hortaliza"

# New row 416: Primera quality record dated 2022-12-05 (serial 44900)
$ws.Range("A416").Value2 = 3
$ws.Range("B416").Value2 = "Femacal de La Calera"
$ws.Range("C416").Value2 = "Coquimbo"
$ws.Range("D416").Value2 = 44900
$ws.Range("E416").Value2 = 5
$ws.Range("F416").Value2 = 100112028
$ws.Range("G416").Value2 = "Sandia"
$ws.Range("H416").Value2 = "Sin especificar"
$ws.Range("I416").Value2 = "Primera"
$ws.Range("J416").Value2 = 520
$ws.Range("K416").Value2 = 3000
$ws.Range("L416").Value2 = 3000
$ws.Range("M416").Value2 = 3000
$ws.Range("N416").Value2 = "`$/unidad"
$ws.Range("O416").Value2 = "Paine"
$ws.Range("P416").Value2 = 3000
$ws.Range("Q416").Value2 = 1
$ws.Range("R416").Value2 = "Hortaliza"

# New row 417: Segunda quality record dated 2022-12-05 (serial 44900)
$ws.Range("A417").Value2 = 3
$ws.Range("B417").Value2 = "Femacal de La Calera"
$ws.Range("C417").Value2 = "Coquimbo"
$ws.Range("D417").Value2 = 44900
$ws.Range("E417").Value2 = 5
$ws.Range("F417").Value2 = 100112028
$ws.Range("G417").Value2 = "Sandia"
$ws.Range("H417").Value2 = "Sin especificar"
$ws.Range("I417").Value2 = "Segunda"
$ws.Range("J417").Value2 = 420
$ws.Range("K417").Value2 = 2000
$ws.Range("L417").Value2 = 2000
$ws.Range("M417").Value2 = 2000
$ws.Range("N417").Value2 = "`$/unidad"
$ws.Range("O417").Value2 = "Paine"
$ws.Range("P417").Value2 = 2000
$ws.Range("Q417").Value2 = 1
$ws.Range("R417").Value2 = "Hortaliza"
